$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.239.01"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "3.491.52"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'604.51"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("D6").Value = "'144.18"
$ws.Range("E6").Value = "  -2.30%  "
$ws.Range("D7").Value = "3.490.36"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").Value = "'8.07"
$ws.Range("E10").Value = "  +2.26%  "
$ws.Range("E11").Value = "  -4.38%  "
$ws.Range("E12").Value = "  -2.29%  "
$ws.Range("D13").Value = "4.084.34"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("E14").Value = "  -4.16%  "
$ws.Range("D15").Value = "'30.37"
$ws.Range("E15").Value = "  -2.50%  "
$ws.Range("D16").Value = "3.486.45"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").Value = "66.253.30"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").Value = "'10.72"
$ws.Range("E19").Value = "  +2.78%  "
$ws.Range("E20").Value = "  -3.43%  "
$ws.Range("E21").Value = "  -3.08%  "
$ws.Range("D22").Value = "'426.26"
$ws.Range("E22").Value = "  -1.69%  "
$ws.Range("D23").Value = "'0.593"
$ws.Range("E23").Value = "  -2.42%  "
$ws.Range("D24").Value = "'77.96"
$ws.Range("E24").Value = "  -2.08%  "
$ws.Range("D25").Value = "3.629.27"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "'0.0000117"
$ws.Range("E27").Value = "  -1.51%  "
$ws.Range("D28").Value = "'9.29"
$ws.Range("E28").Value = "  -5.54%  "
$ws.Range("D29").Value = "'7.93"
$ws.Range("E29").Value = "  -3.89%  "
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("E33").Value = "  -8.39%  "
$ws.Range("D34").Value = "'25.09"
$ws.Range("E34").Value = "  -1.04%  "
$ws.Range("D35").Value = "3.479.52"
$ws.Range("E35").Value = "  -0.34%  "
$ws.Range("D37").Value = "'1.73"
$ws.Range("E37").Value = "  -3.27%  "
$ws.Range("D38").Value = "'5.62"
$ws.Range("E38").Value = "  -4.77%  "
$ws.Range("D39").Value = "'7.73"
$ws.Range("E39").Value = "  -3.40%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").Value = "'170.02"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("D42").Value = "'0.0861"
$ws.Range("E42").Value = "  -3.54%  "
$ws.Range("D43").Value = "'5.16"
$ws.Range("E43").Value = "  -4.71%  "
$ws.Range("D44").Value = "'0.880"
$ws.Range("E44").Value = "  -1.73%  "
$ws.Range("E45").Value = "  -9.00%  "
$ws.Range("D46").Value = "'45.42"
$ws.Range("E46").Value = "  -0.90%  "
$ws.Range("E47").Value = "  -8.69%  "
$ws.Range("D48").Value = "'1.20"
$ws.Range("E48").Value = "  -7.80%  "
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("D50").Value = "'7.14"
$ws.Range("E50").Value = "  -4.23%  "
$ws.Range("D51").Value = "'0.943"
$ws.Range("E51").Value = "  -3.00%  "
